# Update column F ("dSF") values on the active worksheet to reflect
# repulled/recalculated data (repull data, push all data, mean calculation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 1
    5  = 0
    7  = -3
    8  = -1
    10 = -1
    11 = -2
    12 = -3
    13 = -2
    14 = -1
    15 = -2
    16 = -4
    17 = -1
    18 = -6
    19 = -3
    20 = -1
    22 = 1
    23 = -1
    24 = 1
    25 = -1
    26 = -2
    27 = 4
    28 = 2
    29 = 1
    30 = -4
    32 = -1
    33 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
